$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 - shifts existing rows 4..21 down to 5..22
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44742
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112043
$ws.Cells.Item(4, 7).Value = "Pepino dulce"
$ws.Cells.Item(4, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(4, 9).Value = "Segunda"
$ws.Cells.Item(4, 10).Value = 250
$ws.Cells.Item(4, 11).Value = 15000
$ws.Cells.Item(4, 12).Value = 16000
$ws.Cells.Item(4, 13).Value = 15500
$ws.Cells.Item(4, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 861
$ws.Cells.Item(4, 17).Value = 18
$ws.Cells.Item(4, 18).Value = "Hortaliza"
